$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Rows 55-65 on sheet1 ("テスト仕様書"): fill in new use case rows
# Column A: sequence formula; B/C/D/E/F/G: text; I: "?" marker

# Row 55
$ws1.Range('A55').Formula = '=MAX($A$5:A54)+1'
$ws1.Range('B55').Value = '売り切れ商品でチャットできない'
$ws1.Range('C55').Value = '商品詳細/商品一覧のメッセージボタン'
$ws1.Range('D55').Value = '異常系'
$ws1.Range('E55').Value = '売り切れ'
$ws1.Range('F55').Value = '売り切れ商品のメッセージボタンを押下'
$ws1.Range('G55').Value = '売り切れ表示が出てチャット画面へ遷移しない'
$ws1.Range('I55').Value = '?'

# Row 56
$ws1.Range('A56').Formula = '=MAX($A$5:A55)+1'
$ws1.Range('B56').Value = 'レンタル中から購入する'
$ws1.Range('C56').Value = 'レンタル管理（私の申請）'
$ws1.Range('D56').Value = '正常系'
$ws1.Range('E56').Value = '購入確認'
$ws1.Range('F56').Value = 'レンタル中の申請で購入ボタンを押下'
$ws1.Range('G56').Value = '購入金額・レンタル期間の料金・差額の確認画面が表示される'
$ws1.Range('I56').Value = '?'

# Row 57
$ws1.Range('A57').Formula = '=MAX($A$5:A56)+1'
$ws1.Range('D57').Value = '正常系'
$ws1.Range('E57').Value = '購入申請作成'
$ws1.Range('F57').Value = '確認画面で購入を確定'
$ws1.Range('G57').Value = '購入申請が作成され購入管理に表示される（支払額は差額）'
$ws1.Range('I57').Value = '?'

# Row 58
$ws1.Range('A58').Formula = '=MAX($A$5:A57)+1'
$ws1.Range('D58').Value = '正常系'
$ws1.Range('E58').Value = '承認で完了'
$ws1.Range('F58').Value = '出品者がレンタル購入申請を承認'
$ws1.Range('G58').Value = '配送不要で購入完了となりレンタルが完了扱いになる'
$ws1.Range('I58').Value = '?'

# Row 59
$ws1.Range('A59').Formula = '=MAX($A$5:A58)+1'
$ws1.Range('D59').Value = '正常系'
$ws1.Range('E59').Value = '完了表示'
$ws1.Range('F59').Value = '購入完了後にレンタル管理を表示'
$ws1.Range('G59').Value = '購入手続き完了と表示され、配送/購入ボタンが表示されない'
$ws1.Range('I59').Value = '?'

# Row 60
$ws1.Range('A60').Formula = '=MAX($A$5:A59)+1'
$ws1.Range('B60').Value = '購入管理（受け取った注文）'
$ws1.Range('C60').Value = '受注一覧で操作'
$ws1.Range('D60').Value = '正常系'
$ws1.Range('E60').Value = 'レンタル購入承認'
$ws1.Range('F60').Value = 'レンタル中購入の申請を承認'
$ws1.Range('G60').Value = '承認後に配送不要で購入完了になる'
$ws1.Range('I60').Value = '?'

# Row 61
$ws1.Range('A61').Formula = '=MAX($A$5:A60)+1'
$ws1.Range('B61').Value = '取引完了の非表示（受け取った申請）'
$ws1.Range('C61').Value = 'レンタル管理（受け取った申請）'
$ws1.Range('D61').Value = '正常系'
$ws1.Range('E61').Value = '非表示'
$ws1.Range('F61').Value = '完了済み申請の非表示ボタンを押下'
$ws1.Range('G61').Value = '一覧から非表示になり、取引履歴から確認できる'
$ws1.Range('I61').Value = '?'

# Row 62
$ws1.Range('A62').Formula = '=MAX($A$5:A61)+1'
$ws1.Range('B62').Value = '取引完了の非表示（私の申請）'
$ws1.Range('C62').Value = 'レンタル管理（私の申請）'
$ws1.Range('D62').Value = '正常系'
$ws1.Range('E62').Value = '非表示'
$ws1.Range('F62').Value = '完了済み申請の非表示ボタンを押下'
$ws1.Range('G62').Value = '一覧から非表示になり、取引履歴から確認できる'
$ws1.Range('I62').Value = '?'

# Row 63
$ws1.Range('A63').Formula = '=MAX($A$5:A62)+1'
$ws1.Range('B63').Value = '取引完了の非表示（購入管理）'
$ws1.Range('C63').Value = '購入管理（私の購入/受け取った注文）'
$ws1.Range('D63').Value = '正常系'
$ws1.Range('E63').Value = '非表示'
$ws1.Range('F63').Value = '完了済みの購入に非表示ボタンを押下'
$ws1.Range('G63').Value = '一覧から非表示になり、取引履歴から確認できる'
$ws1.Range('I63').Value = '?'

# Row 64
$ws1.Range('A64').Formula = '=MAX($A$5:A63)+1'
$ws1.Range('B64').Value = '取引完了の非表示（返品管理）'
$ws1.Range('C64').Value = '返品管理'
$ws1.Range('D64').Value = '正常系'
$ws1.Range('E64').Value = '非表示'
$ws1.Range('F64').Value = '返品完了/却下の申請に非表示ボタンを押下'
$ws1.Range('G64').Value = '一覧から非表示になり、取引履歴から確認できる'
$ws1.Range('I64').Value = '?'

# Row 65
$ws1.Range('A65').Formula = '=MAX($A$5:A64)+1'
$ws1.Range('B65').Value = '会社概要を表示する'
$ws1.Range('C65').Value = 'サイドバーの会社概要'
$ws1.Range('D65').Value = '正常系'
$ws1.Range('E65').Value = '表示'
$ws1.Range('F65').Value = '会社概要を開く'
$ws1.Range('G65').Value = '会社概要ページが表示される'
$ws1.Range('I65').Value = '?'

# Sheet2 ("エビデンス_"): O1 revision date cell
$ws2.Range('O1').NumberFormat = 'yyyy-mm-dd'
$ws2.Range('O1').Value = 46044
